$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the date stamp in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- "MCF" sheet: bump capacity-factor inputs from 0.95 to 1 ---
$mcf = $wb.Worksheets.Item("MCF")

$cellsToUpdate = @("B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($addr in $cellsToUpdate) {
    $mcf.Range($addr).Value = 1
}

# Move the active selection/cursor to B17 on the MCF sheet (matches last saved view)
$mcf.Activate()
$mcf.Range("B17").Select()
